$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Extend the "times table" grid in rows 22-28 from column AD out to
#    column BF (value = 7*(col-1) - (row-22), same progression already used
#    on the rest of the sheet).
# ---------------------------------------------------------------------------
for ($i = 22; $i -le 28; $i++) {
    for ($j = 31; $j -le 58; $j++) {
        $ws.Cells.Item($i, $j).Value = 7 * ($j - 1) - ($i - 22)
    }
}

# ---------------------------------------------------------------------------
# 2) Re-apply the "highlighted" (prime) cell style to the specific cells
#    that carry it in the target sheet. New cells default to the plain
#    style already, so only these need the highlighted format copied over
#    (copy from E1, an existing highlighted cell, to reuse the same style).
# ---------------------------------------------------------------------------
$highlighted = @(
    @(23,33),
    @(23,37),
    @(23,45),
    @(23,51),
    @(24,34),
    @(24,38),
    @(24,40),
    @(24,46),
    @(24,58),
    @(25,39),
    @(25,41),
    @(25,51),
    @(25,57),
    @(26,34),
    @(26,36),
    @(26,40),
    @(26,42),
    @(26,46),
    @(26,52),
    @(26,54),
    @(27,47),
    @(27,49),
    @(27,53),
    @(27,55),
    @(28,32),
    @(28,42),
    @(28,56)
)

$ws.Range("E1").Copy()
foreach ($p in $highlighted) {
    $ws.Cells.Item($p[0], $p[1]).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Make the used range stretch to BF64 (touch the far corner then clear
#    it again, mirroring how the extent is retained after edits/deletes).
# ---------------------------------------------------------------------------
$ws.Cells.Item(64, 58).Value = "x"
$ws.Cells.Item(64, 58).ClearContents()

# ---------------------------------------------------------------------------
# 4) Update the view: scroll so column J is the left-most visible column and
#    select BG32 as the active cell.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 1
$ws.Cells.Item(32, 59).Select()
